$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 124, shifting existing rows 124:237 down to 125:238
$ws.Rows("124:124").Insert(-4121)

# Populate the newly inserted row 124 with the new record
$ws.Cells.Item(124, 1).Value = 5
$ws.Cells.Item(124, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(124, 3).Value = "Maule"
$ws.Cells.Item(124, 4).Value = 44587
$ws.Cells.Item(124, 5).Value = 7
$ws.Cells.Item(124, 6).Value = 100112003
$ws.Cells.Item(124, 7).Value = "Ajo"
$ws.Cells.Item(124, 8).Value = "Chino"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 200
$ws.Cells.Item(124, 11).Value = 18000
$ws.Cells.Item(124, 12).Value = 18000
$ws.Cells.Item(124, 13).Value = 18000
$ws.Cells.Item(124, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(124, 15).Value = "China"
$ws.Cells.Item(124, 16).Value = 1800
$ws.Cells.Item(124, 17).Value = 10
$ws.Cells.Item(124, 18).Value = "Hortaliza"
